$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would otherwise be auto-parsed as numbers by Excel
# are forced to keep a Text format before the value is written, matching the
# source data which stores every Price/Volume cell as text.

$ws.Range("D2").Value = "67.884.18"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "3.788.03"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.50"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.15"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").Value = "3.786.03"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.07"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "4.420.24"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "3.846.49"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "67.986.33"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.86"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.25"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.698"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000154"
$ws.Range("E24").Value = "  +5.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.58"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -5.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.88"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.21"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "3.739.61"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.993"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.64"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.299"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.36"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.93"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "387.35"
$ws.Range("E50").Value = "  -5.98%  "
$ws.Range("D51").Value = "2.760.50"
$ws.Range("E51").Value = "  +3.49%  "
